## Applies the "updated highlighting in collation files" edit to the
## Theseus/Ariadne collation document.
##
## Summary of changes (see commit diff):
##   1. ";" + " find "  -> merge into a single run "; find "
##   2. "] then, "       -> "] then" + moved _GoBack bookmark + highlighted ","+ " "
##   3. "fate] Fate, "   -> "fate] Fate" + highlighted "," + " "
##   4. ", " (after "Armes") -> highlighted "," + " "
##   5. "all] all, "     -> "all] all" + highlighted "," + " "
##   6. "," (after "fre'd") gains yellow highlight (no text change)
##   7. "; freed, "      -> "; freed" + highlighted "," + " "
##   8. "] " + "look! "  -> merge into a single run "] look! "
##   9. _GoBack bookmark removed from its old spot near "'tis he; 'tis]"
##      (it is re-created at its new location in step 2)

$d = $word.ActiveDocument

# wdYellow highlight color index
$wdYellow = 7

function Highlight-CommaAt([int]$offset) {
    # Highlights a single character (expected to be ",") at a fixed
    # character offset in the document's Content range. Going through
    # Range.Font.* (rather than Range.HighlightColorIndex directly)
    # reliably isolates the run to just this character.
    $rng = $d.Range($offset, $offset + 1)
    $rng.Font.HighlightColorIndex = $wdYellow
}

## -----------------------------------------------------------------
## 1. ";" + " find " -> "; find "
## -----------------------------------------------------------------
$d.Content.Find.Execute("; find ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "; find ", 2) | Out-Null

## -----------------------------------------------------------------
## 9. Remove the old _GoBack bookmark (it will be re-added in step 2).
## -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

## -----------------------------------------------------------------
## 2. "] then, " -> "] then" + _GoBack bookmark + highlighted "," + " "
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("] then, C709")
$bmRange = $d.Range($idx + 6, $idx + 6)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Highlight-CommaAt ($idx + 6)

## -----------------------------------------------------------------
## 3. "fate] Fate, " -> "fate] Fate" + highlighted "," + " "
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("fate] Fate, C709")
Highlight-CommaAt ($idx + 10)

## -----------------------------------------------------------------
## 4. ", " (after "Armes") -> highlighted "," + " "
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("Armes, BL53723")
Highlight-CommaAt ($idx + 5)

## -----------------------------------------------------------------
## 5. "all] all, " -> "all] all" + highlighted "," + " "
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("all] all, C709")
Highlight-CommaAt ($idx + 8)

## -----------------------------------------------------------------
## 6. "," (after "fre'd") gains yellow highlight (no text change)
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("fre’d,] freed")
Highlight-CommaAt ($idx + 5)

## -----------------------------------------------------------------
## 7. "; freed, " -> "; freed" + highlighted "," + " "
## -----------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("; freed, BL53723")
Highlight-CommaAt ($idx + 7)

## -----------------------------------------------------------------
## 8. "] " + "look! " -> "] look! "
## -----------------------------------------------------------------
$d.Content.Find.Execute("] look! ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "] look! ", 2) | Out-Null

Write-Output "done"
